# "Added monte carlo emulation"
#
# The underlying data change is a normalisation of the free-text
# `is_blocked_by` values on the "tasks" sheet so a later (monte-carlo)
# parser can split/trim them reliably:
#   - C5  "TASK-1"        -> "TASK-1  "      (trailing padding spaces)
#   - C8  "TASK-2,TASK-3" -> "TASK-2, TASK-3" (space added after the comma)
#
# Everything else in the xlsx diff (shared-string index shuffling, the
# uniqueCount bump, etc.) is a mechanical side effect of the shared-string
# table being rewritten once these two new unique strings appear - it is
# not an independent edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tasks")

# Order matters: the new unique string for C8 must land in the shared
# string table before the new unique string for C5 to reproduce the
# author's exact table ordering.
$ws.Range("C8").Value = "TASK-2, TASK-3"
$ws.Range("C5").Value = "TASK-1  "

# Cosmetic view-state changes captured in the diff.
$excel.ActiveWindow.Top = 1100
$ws.Range("C16").Select()
